{"js": "// Update the grocery table's Unit Price / Quantity / Total Price columns\n// for each data row (rows 1..6; row 0 is the header).\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index (1-based counting the header as row 0) -> [unitPrice, quantity, totalPrice]\nconst updates = {\n  1: [\"$ 0.97\", \"1\", \"$ 0.97\"], // Banana\n  2: [\"$ 0.72\", \"5\", \"$ 3.60\"], // Strawberry\n  3: [\"$ 0.42\", \"4\", \"$ 1.68\"], // Chicken\n  4: [\"$ 0.05\", \"6\", \"$ 0.30\"], // Bread\n  5: [\"$ 0.77\", \"5\", \"$ 3.85\"], // Eggs\n  6: [\"$ 0.90\", \"2\", \"$ 1.80\"], // Salad\n};\n\nfor (const rowIndexStr of Object.keys(updates)) {\n  const rowIndex = Number(rowIndexStr);\n  const [unitPrice, quantity, totalPrice] = updates[rowIndex];\n  table.getCell(rowIndex, 1).value = unitPrice;\n  table.getCell(rowIndex, 2).value = quantity;\n  table.getCell(rowIndex, 3).value = totalPrice;\n}\n\nawait context.sync();\n", "ps1": "# Update the grocery table's Unit Price / Quantity / Total Price columns\n# for each data row. Row 1 is the header row; data rows are 2..7.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# Word COM table rows/columns are 1-based.\n# Columns: 1 = Product, 2 = Unit Price, 3 = Quantity, 4 = Total Price\n$updates = @{\n    2 = @(\"$ 0.97\", \"1\", \"$ 0.97\")  # Banana\n    3 = @(\"$ 0.72\", \"5\", \"$ 3.60\")  # Strawberry\n    4 = @(\"$ 0.42\", \"4\", \"$ 1.68\")  # Chicken\n    5 = @(\"$ 0.05\", \"6\", \"$ 0.30\")  # Bread\n    6 = @(\"$ 0.77\", \"5\", \"$ 3.85\")  # Eggs\n    7 = @(\"$ 0.90\", \"2\", \"$ 1.80\")  # Salad\n}\n\nforeach ($row in $updates.Keys) {\n    $values = $updates[$row]\n    $t.Cell($row, 2).Range.Text = $values[0]\n    $t.Cell($row, 3).Range.Text = $values[1]\n    $t.Cell($row, 4).Range.Text = $values[2]\n}\n"}
